# Estado de cuenta - actualiza base de datos EC y agrega parte 1 de nuevo periodo (2508)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row right below the existing data row (row 16), shifting the
#    signature block (old rows 21/22) down to 22/23.
$ws.Rows.Item(17).Insert()

# 2. Seed the new row with the same look & feel (borders, number formats, font)
#    and worker identity data as the row above, then overwrite the period /
#    value columns for the new period.
$srcRow = $ws.Range("B16:J16")
$dstRow = $ws.Range("B17:J17")
$srcRow.Copy($dstRow)

$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# 3. Existing period's "Valor Mora" changed.
$ws.Range("G16").Value = 1423500

# 4. Header totals: total "VALOR MORA" and "Cant. Periodos" both reflect the
#    newly added period.
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2
